# Exclude denied embassy requests
# The region "The beforelife" (row 68) had its embassy request denied, so
# remove that entire row and let the remaining rows shift up.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("A68:I68").EntireRow.Delete()
